$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update order fraction counters (merged cells H7:K7 and H8:K8 store value in top-left cell)
$ws.Range("H7").Value = "1:0"
$ws.Range("H8").Value = "1:2"

# Update the generated-at timestamp (merged cell A10:F10)
$ws.Range("A10").Value = "Monday, 15 September, 2025 9:33 AM"
